# Updated to fetch data from excel using DataProvider
#
# - Renames "Deals" -> "new_deal" and trims its stray padding columns
#   (G:Q) down to A:F, adding a new "Exclude Reports?" column (text "true").
# - Adds a brand-new "edit_deal" sheet (between new_deal and
#   login_with_invalid_data) with the same header row plus an
#   "Exclude Reports?" boolean column and the old "Updated Title" /
#   "Updated deal" columns.
# - Leaves "login_with_invalid_data" untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Deals" -> "new_deal"
# ---------------------------------------------------------------------
$dealsSheet = $wb.Worksheets.Item("Deals")
$dealsSheet.Name = "new_deal"

# Clear the old stray styled-but-empty padding cells (G1:Q2) so the sheet
# shrinks back down to A1:F2.
$dealsSheet.Range("G1:Q2").Clear()

# New header/data cell for the "Exclude Reports?" column.
$dealsSheet.Range("F1").Value = "Exclude Reports?"
# Force literal text (not boolean) via the classic apostrophe-prefix.
$dealsSheet.Range("F2").Value = "'true"

# Re-apply the header / data styling (yellow+border header row, bordered
# data row) plus a text number format on the touched cells, matching the
# rest of the sheet's look.
$dealsSheet.Range("A1:F2").NumberFormat = "@"

$dealsSheet.Range("A1:F1").Interior.ColorIndex = 6
$dealsSheet.Range("A1:F2").Borders.LineStyle = 1

$dealsSheet.Range("B26").Select()

# ---------------------------------------------------------------------
# 2) New "edit_deal" sheet, inserted between "new_deal" and
#    "login_with_invalid_data"
# ---------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("login_with_invalid_data")
$editDeal = $wb.Worksheets.Add($loginSheet)
$editDeal.Name = "edit_deal"

$editDeal.Range("A1").Value = "Title"
$editDeal.Range("B1").Value = "Amount"
$editDeal.Range("C1").Value = "Quantity"
$editDeal.Range("D1").Value = "Type"
$editDeal.Range("E1").Value = "Status"
$editDeal.Range("F1").Value = "Exclude Reports?"
$editDeal.Range("G1").Value = "Updated Title"

$editDeal.Range("A2").Value = "Deal Number one"
$editDeal.Range("B2").Value = 2000
$editDeal.Range("C2").Value = 2
$editDeal.Range("D2").Value = "Old"
$editDeal.Range("E2").Value = "Closed"
$editDeal.Range("F2").Value = $true
$editDeal.Range("G2").Value = "Updated deal"

$editDeal.Range("A1:G2").NumberFormat = "@"
$editDeal.Range("A1:G1").Interior.ColorIndex = 6
$editDeal.Range("A1:G2").Borders.LineStyle = 1

$editDeal.Range("F2").Select()

# Columns widths to roughly match the authored layout.
$editDeal.Columns.Item(1).ColumnWidth = 13.42578125
$editDeal.Columns.Item(6).ColumnWidth = 19.140625
$editDeal.Columns.Item(7).ColumnWidth = 13.5703125

$dealsSheet.Columns.Item(1).ColumnWidth = 13.42578125
$dealsSheet.Columns.Item(6).ColumnWidth = 19.140625

# Make "edit_deal" the active tab, matching activeTab="1" in the workbook.
$editDeal.Activate()
